$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 15 (A15) down to row 16 (A16) so the new row
# matches the existing style (bold/centered/bordered first column).
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 0.9970734475936542
$ws.Range("D16").Value = 1.001684715932086
$ws.Range("E16").Value = 0.9964559641710516
$ws.Range("F16").Value = 0.9970734475936542
$ws.Range("G16").Value = 0.9982352941176471
$ws.Range("H16").Value = 0.9964413401068091
$ws.Range("I16").Value = 0.9952941176470588
$ws.Range("J16").Value = 1.001684715932086
$ws.Range("K16").Value = 0.9990703400515688
$ws.Range("L16").Value = 0.9980718938226114
$ws.Range("M16").Value = 0.9975308132613844
